$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "11/06/2025"
$ws.Range("A66").ClearFormats()
$ws.Range("B66").Value = 0.2107717439314687
$ws.Range("C66").Value = 0.7892282560685313
